$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first two data rows (old rows 2 and 3), shifting all subsequent
# rows up by two. This turns old row 4 into new row 2, ... old row 22 into
# new row 20, and shrinks the used range from A1:H22 to A1:H21.
$ws.Range("A2:H3").EntireRow.Delete() | Out-Null

# The new last row (row 21) receives freshly generated sensor data rather
# than reusing any previously shifted values; timestamp/label are re-applied
# since the row-delete above left this trailing row blank.
$ws.Range("A21").Value = 1900
$ws.Range("B21").Value = "falling"
$ws.Range("C21").Value = 0.1471566200256338
$ws.Range("D21").Value = 1.182808732986451
$ws.Range("E21").Value = 0.03668105900287391
$ws.Range("F21").Value = 0.0018325957935303
$ws.Range("G21").Value = 0.0178678091615438
$ws.Range("H21").Value = 0.0360410511493682
